$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(375, 44449, 0, 4, 47.13090609166961),
    @(376, 44450, 0, 4, 47.13090609166961),
    @(377, 44451, 2, 6, 70.69635913750442),
    @(378, 44452, 1, 7, 82.47908566042182),
    @(379, 44453, 0, 6, 70.69635913750442),
    @(380, 44454, 0, 6, 70.69635913750442),
    @(381, 44455, 1, 4, 47.13090609166961),
    @(382, 44456, 0, 4, 47.13090609166961),
    @(383, 44457, 1, 5, 58.91363261458702),
    @(384, 44458, 1, 4, 47.13090609166961),
    @(385, 44459, 0, 3, 35.34817956875221)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Carry the date-column formatting (style index 2 used by A2:A374) down
# onto the newly-appended rows, matching column A of the preceding row.
$ws.Range("A374").Copy()
$ws.Range("A375:A385").PasteSpecial(-4122)
